$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Bump the "datetimeFigureOut" date placeholder shown on every slide
#    layout (and on the slide master) from 4/15/2019 to 4/25/2019.
# ---------------------------------------------------------------------------

# Slide master's own "Date Placeholder" shape.
$masterShapes = $p.SlideMaster.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $sh = $masterShapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "4/15/2019") {
                $sh.TextFrame.TextRange.Text = "4/25/2019"
            }
        }
    }
}

# Every slide layout hanging off the master has its own "Date Placeholder".
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame) {
                if ($sh.TextFrame.TextRange.Text -eq "4/15/2019") {
                    $sh.TextFrame.TextRange.Text = "4/25/2019"
                }
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 7 ("Thank you for completing the practice stage ...") - update the
#    task duration from 5 minutes to 25 minutes.
# ---------------------------------------------------------------------------

$slide = $p.Slides.Item(7)
$shape = $slide.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange
$full = $tr.Text
$needle = "take 5 "
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, $needle.Length)
    $sub.Text = "take 25 "
}
